$d = $word.ActiveDocument

# The "_GoBack" bookmark currently sits at the very end of the paragraph whose
# text is "12345". The edit splits that paragraph right after "12345": a new
# paragraph is inserted, containing the new run "123456", and the bookmark
# (both bookmarkStart and bookmarkEnd) moves to trail that new run.

$bm = $d.Bookmarks.Item("_GoBack")
$bmPara = $bm.Range.Paragraphs(1)

# Insertion point = end of the paragraph's text, i.e. just before its
# paragraph mark (Range.End includes the paragraph mark, so back up by one).
$paraEnd = $bmPara.Range.End
$textEnd = $paraEnd - 1
$splitPoint = $d.Range($textEnd, $textEnd)

# Split the paragraph: this creates a new, empty paragraph right after the
# "12345" paragraph and leaves the bookmark in place (still right after
# "12345") for now.
$splitPoint.InsertParagraphAfter()

$newParaIndex = $bmPara.Index + 1
$newPara = $d.Paragraphs($newParaIndex)

# Populate the new paragraph with "123456" plus a trailing sentinel
# character; the sentinel gives us a non-boundary (mid-paragraph) position
# to (re)plant the "_GoBack" bookmark around, which this host handles more
# reliably than a zero-length bookmark sitting exactly at a paragraph's text
# end. Inserting this way also makes the new run inherit the same rPr
# (lang="en-US") as its neighbours.
$newPara.Range.InsertBefore("123456X")

$newPara = $d.Paragraphs($newParaIndex)
$newParaEnd = $newPara.Range.End
$sentinelStart = $newParaEnd - 2
$sentinelEnd = $newParaEnd - 1
$sentinelRange = $d.Range($sentinelStart, $sentinelEnd)

# Re-anchor "_GoBack" around the sentinel character -- Bookmarks.Add with an
# existing bookmark name relocates it rather than creating a duplicate.
$null = $d.Bookmarks.Add("_GoBack", $sentinelRange)

# Remove the sentinel character; the now-collapsed bookmark is left sitting
# right after "123456", matching the target layout.
$cleanupRange = $d.Range($sentinelStart, $sentinelEnd)
$cleanupRange.Delete()
